$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the "Ins"/"Del" labels to "I"/"D". The first two (A2, A3) are
# re-entered with a leading apostrophe so Excel stores them as
# quote-prefixed text (quotePrefix="1"), matching the target formatting.
$ws.Range("A2").Value = "'I 100"
$ws.Range("A3").Value = "'I 1,000"
$ws.Range("A4").Value = "I 10,000"
$ws.Range("A5").Value = "I 100,000"
$ws.Range("A6").Value = "D 100"
$ws.Range("A7").Value = "D 1,000"
$ws.Range("A8").Value = "D 10,000"
$ws.Range("A9").Value = "D 100,000"

# Move/restore the active selection as it was left in the source file.
[void]$ws.Range("C16").Select()
